# Update the "想去人数" (want-to-go count) figures that changed between
# the two published snapshots of the data. The same table is duplicated
# on the "展览" sheet and the "全部类型" sheet, so both need updating.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of cell address -> new value for each of the two sheets.
$updates = @{
    "F4"  = 12158
    "F11" = 439
    "F17" = 2626
    "F18" = 88
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
